$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.442.60"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.891.85"
$ws.Range("E3").Value = "  +0.02%  "
$dStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = $dStyle
$ws.Range("E4").Value = "  -0.02%  "
$dStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.81"
$ws.Range("D5").Style = $dStyle
$ws.Range("E5").Value = "  -0.28%  "
$dStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = $dStyle
$ws.Range("E6").Value = "  -0.01%  "
$dStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4847"
$ws.Range("D7").Style = $dStyle
$ws.Range("E7").Value = "  -1.27%  "
$dStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("D8").Style = $dStyle
$ws.Range("E8").Value = "  -1.46%  "
$dStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06614"
$ws.Range("D9").Style = $dStyle
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "1.897.43"
$ws.Range("E10").Value = "  -0.17%  "
$dStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.86"
$ws.Range("D11").Style = $dStyle
$ws.Range("E11").Value = "  -0.95%  "
$dStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07402"
$ws.Range("D12").Style = $dStyle
$ws.Range("E12").Value = "  +0.81%  "
$dStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.200"
$ws.Range("D13").Style = $dStyle
$ws.Range("E13").Value = "  +0.61%  "
$dStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.85"
$ws.Range("D14").Style = $dStyle
$ws.Range("E14").Value = "  +0.93%  "
$dStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6625"
$ws.Range("D15").Style = $dStyle
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "30.418.02"
$ws.Range("E16").Value = "  -0.55%  "
$dStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.55"
$ws.Range("D17").Style = $dStyle
$ws.Range("E17").Value = "  +0.78%  "
$dStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007784"
$ws.Range("D18").Style = $dStyle
$ws.Range("E18").Value = "  -1.20%  "
$dStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = $dStyle
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "2.138.05"
$ws.Range("E20").Value = "  -0.69%  "
$dStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.379"
$ws.Range("D21").Style = $dStyle
$ws.Range("E21").Value = "  +0.34%  "
$dStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = $dStyle
$ws.Range("E22").Value = "  +0.04%  "
$dStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "220.55"
$ws.Range("D23").Style = $dStyle
$ws.Range("E23").Value = "  +16.81%  "
$dStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.240"
$ws.Range("D24").Style = $dStyle
$ws.Range("E24").Value = "  +0.52%  "
$dStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.392"
$ws.Range("D25").Style = $dStyle
$ws.Range("E25").Value = "  -1.51%  "
$dStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.98"
$ws.Range("D26").Style = $dStyle
$ws.Range("E26").Value = "  +0.78%  "
$dStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.22"
$ws.Range("D27").Style = $dStyle
$ws.Range("E27").Value = "  +4.04%  "
$dStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.946"
$ws.Range("D28").Style = $dStyle
$ws.Range("E28").Value = "  +0.77%  "
$dStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.452"
$ws.Range("D29").Style = $dStyle
$ws.Range("E29").Value = "  -0.95%  "
$dStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.334"
$ws.Range("D30").Style = $dStyle
$ws.Range("E30").Value = "  -1.92%  "
$dStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09219"
$ws.Range("D31").Style = $dStyle
$ws.Range("E31").Value = "  +0.67%  "
$dStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.040"
$ws.Range("D32").Style = $dStyle
$ws.Range("E32").Value = "  -0.27%  "
$dStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05080"
$ws.Range("D33").Style = $dStyle
$ws.Range("E33").Value = "  -3.18%  "
$dStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7597"
$ws.Range("D34").Style = $dStyle
$ws.Range("E34").Value = "  +2.25%  "
$dStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("D35").Style = $dStyle
$ws.Range("E35").Value = "  +4.94%  "
$dStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.702"
$ws.Range("D36").Style = $dStyle
$ws.Range("E36").Value = "  -1.03%  "
$dStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01880"
$ws.Range("D37").Style = $dStyle
$ws.Range("E37").Value = "  +3.02%  "
$dStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.645"
$ws.Range("D38").Style = $dStyle
$ws.Range("E38").Value = "  -2.00%  "
$dStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9199"
$ws.Range("D39").Style = $dStyle
$ws.Range("E39").Value = "  +0.51%  "
$dStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.086"
$ws.Range("D40").Style = $dStyle
$ws.Range("E40").Value = "  +0.89%  "
$dStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.963"
$ws.Range("D41").Style = $dStyle
$ws.Range("E41").Value = "  +0.63%  "
$dStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4364"
$ws.Range("D42").Style = $dStyle
$ws.Range("E42").Value = "  -1.32%  "
$dStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.18"
$ws.Range("D43").Style = $dStyle
$ws.Range("E43").Value = "  -0.02%  "
$dStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = $dStyle
$ws.Range("E44").Value = "  +1.06%  "
$dStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.636"
$ws.Range("D45").Style = $dStyle
$ws.Range("E45").Value = "  +1.03%  "
$dStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.612"
$ws.Range("D46").Style = $dStyle
$ws.Range("E46").Value = "  +13.18%  "
$dStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1328"
$ws.Range("D47").Style = $dStyle
$ws.Range("E47").Value = "  -3.71%  "
$dStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.75"
$ws.Range("D48").Style = $dStyle
$ws.Range("E48").Value = "  -12.60%  "
$dStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.955"
$ws.Range("D49").Style = $dStyle
$ws.Range("E49").Value = "  -0.77%  "
$dStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.63"
$ws.Range("D50").Style = $dStyle
$ws.Range("E50").Value = "  -2.41%  "
$dStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05720"
$ws.Range("D51").Style = $dStyle
$ws.Range("E51").Value = "  -2.09%  "
